$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.936.62"
$ws.Range("E2").Value = "'  -0.36%  "
$ws.Range("D3").Value = "'2.497.14"
$ws.Range("E3").Value = "'  -0.73%  "
$ws.Range("E4").Value = "'  -0.13%  "
$ws.Range("D5").Value = "'535.31"
$ws.Range("E5").Value = "'  +0.22%  "
$ws.Range("D6").Value = "'137.23"
$ws.Range("E6").Value = "'  -0.97%  "
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "'  -0.42%  "
$ws.Range("E8").Value = "'  +1.01%  "
$ws.Range("D9").Value = "'2.518.08"
$ws.Range("E9").Value = "'  -0.03%  "
$ws.Range("E10").Value = "'  +2.26%  "
$ws.Range("E11").Value = "'  -0.22%  "
$ws.Range("D12").Value = "'5.35"
$ws.Range("E12").Value = "'  -0.67%  "
$ws.Range("E13").Value = "'  -1.79%  "
$ws.Range("D14").Value = "'2.943.35"
$ws.Range("E14").Value = "'  -0.77%  "
$ws.Range("D15").Value = "'23.22"
$ws.Range("E15").Value = "'  +1.18%  "
$ws.Range("D16").Value = "'58.870.62"
$ws.Range("E16").Value = "'  -0.39%  "
$ws.Range("D17").Value = "'0.0000140"
$ws.Range("E17").Value = "'  +0.27%  "
$ws.Range("D18").Value = "'2.512.72"
$ws.Range("E18").Value = "'  -1.93%  "
$ws.Range("D19").Value = "'11.06"
$ws.Range("E19").Value = "'  +1.85%  "
$ws.Range("D20").Value = "'4.26"
$ws.Range("E20").Value = "'  +1.77%  "
$ws.Range("D21").Value = "'325.40"
$ws.Range("E21").Value = "'  +1.74%  "
$ws.Range("E22").Value = "'  +0.12%  "
$ws.Range("D23").Value = "'5.88"
$ws.Range("E23").Value = "'  +1.97%  "
$ws.Range("D24").Value = "'65.07"
$ws.Range("E24").Value = "'  +4.68%  "
$ws.Range("D25").Value = "'0.421"
$ws.Range("E25").Value = "'  +0.11%  "
$ws.Range("D26").Value = "'0.166"
$ws.Range("E26").Value = "'  +0.81%  "
$ws.Range("D27").Value = "'0.996"
$ws.Range("E27").Value = "'  -0.50%  "
$ws.Range("D28").Value = "'7.59"
$ws.Range("E28").Value = "'  -1.76%  "
$ws.Range("D29").Value = "'6.76"
$ws.Range("E29").Value = "'  +1.35%  "
$ws.Range("D30").Value = "'0.0₃0776"
$ws.Range("E30").Value = "'  +2.24%  "
$ws.Range("D31").Value = "'1.77"
$ws.Range("E31").Value = "'  -1.22%  "
$ws.Range("D32").Value = "'167.26"
$ws.Range("E32").Value = "'  +4.55%  "
$ws.Range("E33").Value = "'  +5.22%  "
$ws.Range("E34").Value = "'  -0.16%  "
$ws.Range("E35").Value = "'  -2.69%  "
$ws.Range("D36").Value = "'18.58"
$ws.Range("E36").Value = "'  +0.73%  "
$ws.Range("D37").Value = "'4.12"
$ws.Range("E37").Value = "'  -0.81%  "
$ws.Range("D38").Value = "'1.57"
$ws.Range("E38").Value = "'  -0.10%  "
$ws.Range("D39").Value = "'36.80"
$ws.Range("E39").Value = "'  -0.25%  "
$ws.Range("D40").Value = "'0.827"
$ws.Range("E40").Value = "'  +3.43%  "
$ws.Range("D41").Value = "'3.63"
$ws.Range("E41").Value = "'  +0.32%  "
$ws.Range("D42").Value = "'5.27"
$ws.Range("E42").Value = "'  +1.58%  "
$ws.Range("D43").Value = "'280.23"
$ws.Range("E43").Value = "'  -0.61%  "
$ws.Range("D44").Value = "'0.994"
$ws.Range("E44").Value = "'  -0.61%  "
$ws.Range("D45").Value = "'0.605"
$ws.Range("E45").Value = "'  +1.70%  "
$ws.Range("D46").Value = "'10.86"
$ws.Range("E46").Value = "'  -0.09%  "
$ws.Range("D47").Value = "'128.20"
$ws.Range("E47").Value = "'  +4.42%  "
$ws.Range("D48").Value = "'0.0930"
$ws.Range("E48").Value = "'  +1.00%  "
$ws.Range("D49").Value = "'0.0515"
$ws.Range("E49").Value = "'  +1.72%  "
$ws.Range("E50").Value = "'  +0.43%  "
$ws.Range("D51").Value = "'17.40"
$ws.Range("E51").Value = "'  +0.48%  "
